$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Friendlier row labels (column A) and header captions
# ---------------------------------------------------------------------------
$ws.Range("A2").Value  = "Is post-quantum?"
$ws.Range("A3").Value  = "Execution times"
$ws.Range("A4").Value  = "Signature size (B)"
$ws.Range("A5").Value  = "Public key size (B)"
$ws.Range("A6").Value  = "Private key size (B)"
$ws.Range("A7").Value  = "Keys gen time mean (ns)"
$ws.Range("A8").Value  = "Keys gen time standard deviation (ns)"
$ws.Range("A9").Value  = "Signature time mean (ns)"
$ws.Range("A10").Value = "Signature time standard deviation (ns)"
$ws.Range("A11").Value = "Verify time mean (ns)"
$ws.Range("A12").Value = "Verify time standard deviation (ns)"

# ---------------------------------------------------------------------------
# 2) Updated benchmark numbers
# ---------------------------------------------------------------------------
$ws.Range("B3:F3").Value = 100

$ws.Range("E4").Value = 655

$ws.Range("B6").Value = 635

$ws.Range("B7").Value = 28950392
$ws.Range("C7").Value = 968381
$ws.Range("D7").Value = 39135
$ws.Range("E7").Value = 8997841
$ws.Range("F7").Value = 420322

$ws.Range("B8").Value = 16790592
$ws.Range("C8").Value = 144602
$ws.Range("D8").Value = 2864
$ws.Range("E8").Value = 1904224
$ws.Range("F8").Value = 25897

$ws.Range("B9").Value = 390125
$ws.Range("C9").Value = 575886
$ws.Range("D9").Value = 94352
$ws.Range("E9").Value = 367723
$ws.Range("F9").Value = 9761512

$ws.Range("B10").Value = 70077
$ws.Range("C10").Value = 84300
$ws.Range("D10").Value = 50797
$ws.Range("E10").Value = 16987
$ws.Range("F10").Value = 341901

$ws.Range("B11").Value = 76659
$ws.Range("C11").Value = 1866386
$ws.Range("D11").Value = 36509
$ws.Range("E11").Value = 68632
$ws.Range("F11").Value = 867417

$ws.Range("B12").Value = 17937
$ws.Range("C12").Value = 153704
$ws.Range("D12").Value = 2443
$ws.Range("E12").Value = 2715
$ws.Range("F12").Value = 33777

# ---------------------------------------------------------------------------
# 3) Prettier formatting
# ---------------------------------------------------------------------------

# Recolor the header/label fill (indexed 48 -> indexed 22) and turn wrap-text
# on for the label column - this is the style already used by column A.
$ws.Range("A1:A12").Interior.ColorIndex = 22
$ws.Range("A1:A12").WrapText = $true

# Make the top-of-table algorithm names share that same label look by
# copying the finished column-A format onto them.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:F1").PasteSpecial(-4122) | Out-Null

# Data cells: smaller Arial font, wrapped and right-aligned.
$data = $ws.Range("B2:F12")
$data.Font.Name = "Arial"
$data.Font.Size = 12
$data.WrapText = $true
$data.HorizontalAlignment = -4152

# Taller rows so the wrapped text reads comfortably.
$ws.Range("A1:A12").EntireRow.RowHeight = 37.5

$excel.CutCopyMode = $false
